$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 26 (RM 232), then delete the row that becomes row 27 (SC 92),
# which shifts all subsequent rows up by two.
$ws.Rows.Item(26).Delete()
$ws.Rows.Item(27).Delete()

# Apply individual cell-value edits (post row-shift) to match target state.
$ws.Range("D5").Value = $null
$ws.Range("C6").Value = 15.1
$ws.Range("C8").Value = $null
$ws.Range("D11").Value = -15.5
$ws.Range("C19").Value = 13.2
$ws.Range("D19").Value = $null
$ws.Range("C21").Value = $null
$ws.Range("C23").Value = 12.2
$ws.Range("D23").Value = -13.9
$ws.Range("D25").Value = -15.5

$ws.Range("B26").Value = $null
$ws.Range("B27").Value = -20.4
$ws.Range("C27").Value = $null
$ws.Range("D27").Value = $null
$ws.Range("B29").Value = $null
$ws.Range("C29").Value = 11.2
$ws.Range("D29").Value = $null
$ws.Range("D30").Value = -13.6
$ws.Range("D33").Value = -14.1
